$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "K" column (G) values regenerated to use K instead of Strike# (TB),
# recomputed from the regenerated std/mean + s_vals calculation.
$kValues = @(
    1  # G2
    0  # G3
    0  # G4
    1  # G5
    2  # G6
    2  # G7
    0  # G8
    0  # G9
    2  # G10
    2  # G11
    1  # G12
    2  # G13
    1  # G14
    2  # G15
    1  # G16
    0  # G17
    0  # G18
    1  # G19
    1  # G20
    1  # G21
    0  # G22
    0  # G23
    1  # G24
    0  # G25
    2  # G26
    1  # G27
    1  # G28
    2  # G29
    1  # G30
    2  # G31
    0  # G32
    0  # G33
    1  # G34
    0  # G35
    1  # G36
    0  # G37
    1  # G38
    0  # G39
    0  # G40
    0  # G41
    0  # G42
    2  # G43
    0  # G44
    0  # G45
    1  # G46
    0  # G47
    2  # G48
    1  # G49
    1  # G50
    1  # G51
    1  # G52
    0  # G53
    0  # G54
    1  # G55
    0  # G56
    2  # G57
    0  # G58
    1  # G59
    0  # G60
    1  # G61
    1  # G62
    0  # G63
    1  # G64
    1  # G65
    1  # G66
    2  # G67
    1  # G68
    2  # G69
    0  # G70
    0  # G71
)

for ($i = 0; $i -lt $kValues.Count; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 7).Value = $kValues[$i]
}

